# Update "想去人数" (number of people interested) figures in column F
# for rows 2-13 on both the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 416
    3  = 1407
    4  = 7086
    5  = 531
    6  = 289
    7  = 4805
    8  = 119
    9  = 1177
    10 = 51
    11 = 960
    12 = 256
    13 = 5442
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
